$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the new "ArrestWarrant" rule row (row 25) that mirrors the
# existing rule rows above it.
$ws.Range("B25").Value = "ArrestWarrant"
$ws.Range("C25").Value = "file "
$ws.Range("D25").Value = "true "
$ws.Range("E25").Value = "arrestWarrant"
$ws.Range("F25").Value = 50
$ws.Range("G25").Value = "P3D"
$ws.Range("H25").Value = "Review Arrest Warrant"
$ws.Range("I25").Value = "ann-acm@armedia.com,ian-acm@armedia.com,samuel-acm@armedia.com"
$ws.Range("K25").Value = "true "

# Turn the approvers cell into a mailto hyperlink, same as I22.
$ws.Hyperlinks.Add($ws.Range("I25"), "mailto:ann-acm@armedia.com,ian-acm@armedia.com,samuel-acm@armedia.com")

# Reflect the new cursor/scroll position left in the saved view.
$ws.Range("J25").Select()
